# Daily attendance processing - sorts the "Recorded By" list (column G)
# of each data row alphabetically (case-insensitive, with an ordinal
# tie-break so that e.g. "System" sorts before "system").

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$used = $ws.UsedRange
$totalRows = $used.Rows.Count

for ($rowNum = 2; $rowNum -le $totalRows; $rowNum++) {
    $origText = $ws.Cells.Item($rowNum, 7).Text
    if ([string]::IsNullOrEmpty($origText)) { continue }

    $nameParts = $origText -split ", "
    $partCount = $nameParts.Count
    if ($partCount -le 1) { continue }

    # Insertion sort: primary key = lower-case value, tie-break = ordinal
    # (character code) comparison of the original value.
    for ($ii = 1; $ii -lt $partCount; $ii++) {
        $curKey = $nameParts[$ii]
        $curKeyLower = $curKey.ToLower()
        $jj = $ii - 1
        $keepGoing = $true
        while ($keepGoing -eq $true) {
            if ($jj -lt 0) {
                $keepGoing = $false
            } else {
                $cmpVal = $nameParts[$jj]
                $cmpValLower = $cmpVal.ToLower()
                $isGreater = $false
                if ($cmpValLower -gt $curKeyLower) {
                    $isGreater = $true
                } elseif ($cmpValLower -eq $curKeyLower) {
                    $minLen = $cmpVal.Length
                    if ($curKey.Length -lt $minLen) { $minLen = $curKey.Length }
                    $charIdx = 0
                    $decided = $false
                    while ($charIdx -lt $minLen -and $decided -eq $false) {
                        $codeA = [int][char]$cmpVal[$charIdx]
                        $codeB = [int][char]$curKey[$charIdx]
                        if ($codeA -gt $codeB) {
                            $isGreater = $true
                            $decided = $true
                        } elseif ($codeA -lt $codeB) {
                            $isGreater = $false
                            $decided = $true
                        }
                        $charIdx = $charIdx + 1
                    }
                    if ($decided -eq $false -and $cmpVal.Length -gt $curKey.Length) {
                        $isGreater = $true
                    }
                }

                if ($isGreater -eq $true) {
                    $nameParts[$jj + 1] = $nameParts[$jj]
                    $jj = $jj - 1
                } else {
                    $keepGoing = $false
                }
            }
        }
        $nameParts[$jj + 1] = $curKey
    }

    $newText = [string]::Join(", ", $nameParts)
    if ($newText -ne $origText) {
        $ws.Range("G" + $rowNum).Value = $newText
    }
}
